# Generate Report for Handback
#
# - The status for file "5ba90bdc-6db1-4b00-b2b1-7eb0703df7c9.md" changes
#   from "Ready for handoff" to "Handback transform failed" (shared string
#   used by the Overview sheet's row for that file).
# - An "Error Detail" message is recorded for that same file's row on the
#   zh-cn and de-de locale sheets (column L, row 7).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the status text wherever it currently reads "Ready for handoff"
# (the Overview sheet's zh-cn/de-de status columns, plus the Status column
# on each locale sheet's row for the same file).
$overview.Range("B7").Value = "Handback transform failed"
$overview.Range("C7").Value = "Handback transform failed"
$zhcn.Range("C7").Value = "Handback transform failed"
$dede.Range("C7").Value = "Handback transform failed"

# Record the handback/handoff file name mismatch error for each locale.
$zhcn.Range("L7").Value = "Handback file name: g2ro1uby.ugy is different with handoff file name: 5ba90bdc-6db1-4b00-b2b1-7eb0703df7c9.01e2f040f9f62ff0a366b88afb0f8b2bbc0309c2.zh-cn."
$dede.Range("L7").Value = "Handback file name: g2ro1uby.ugy is different with handoff file name: 5ba90bdc-6db1-4b00-b2b1-7eb0703df7c9.01e2f040f9f62ff0a366b88afb0f8b2bbc0309c2.de-de."
